# Insert a new price-observation row for "Albahaca" (Vega Modelo de Temuco)
# right before the current row 230, pushing the existing rows 230-256 down
# to 231-257 and extending the used range to A1:R257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 230 (shifts rows 230..256 down to 231..257)
$ws.Range("A230").EntireRow.Insert()

# Populate the newly inserted row 230 with the new weekly observation
$ws.Range("A230").Value = 10
$ws.Range("B230").Value = "Vega Modelo de Temuco"
$ws.Range("C230").Value = "La Araucanía"
$ws.Range("D230").Value = 44776
$ws.Range("E230").Value = 9
$ws.Range("F230").Value = 100112052
$ws.Range("G230").Value = "Albahaca"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 50
$ws.Range("K230").Value = 6000
$ws.Range("L230").Value = 6000
$ws.Range("M230").Value = 6000
$ws.Range("N230").Value = "`$/paquete"
$ws.Range("O230").Value = "Región de Arica y Parinacota"
$ws.Range("P230").Value = 6000
$ws.Range("Q230").Value = 1
$ws.Range("R230").Value = "Hortaliza"
